$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a second "Formula2" header in B1, reusing A1's formatting (style s="1")
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "Formula2"

# Mirror the numeric series from column A into column B (rows 2-9)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}
